$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price (D) and 1h volume change (E) values
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.794.77'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.854.68'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -1.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.21'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4301'
$ws.Range("E7").Value = '  -1.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3759'
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07351'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8785'
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.60'
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.859.82'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.769'
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.452'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07131'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.97'
$ws.Range("E16").Value = '  +4.48%  '
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009018'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.010'
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.805.95'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.226'
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.08'
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.083.82'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.988'
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.41'
$ws.Range("E26").Value = '  -1.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.64'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("E28").Value = '  +9.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.372'
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.18'
$ws.Range("E30").Value = '  +1.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08946'
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("E32").Value = '  +1.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7792'
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.551'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.933'
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.012'
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.132'
$ws.Range("E37").Value = '  -0.64%  '
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05351'
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("E40").Value = '  +1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.160'
$ws.Range("E41").Value = '  +4.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1698'
$ws.Range("E42").Value = '  +1.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5143'
$ws.Range("E43").Value = '  -0.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.841'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.73'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.96'
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4770'
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06480'
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.694'
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.851'
$ws.Range("E51").Value = '  -2.52%  '
